$d = $word.ActiveDocument

function Escape-Xml {
    param([string]$Text)
    $t = $Text -replace '&', '&amp;'
    $t = $t -replace '<', '&lt;'
    $t = $t -replace '>', '&gt;'
    return $t
}

function Set-ParagraphRuns {
    param([string]$MatchText, [string]$InnerRunsXml)
    foreach ($p in $d.Paragraphs) {
        $paraText = $p.Range.Text.TrimEnd([char]13)
        if ($paraText -eq $MatchText) {
            # Select the paragraph's content (excluding the trailing paragraph mark)
            # and replace it via InsertXML with an explicit run-level fragment so the
            # new content lands inside this exact paragraph instead of spilling into
            # a sibling paragraph.
            $rng = $d.Range($p.Range.Start, $p.Range.End - 1)
            $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
                   '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
                   '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                   '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body><w:p>' + $InnerRunsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
            $rng.InsertXML($pkg) | Out-Null
            return $true
        }
    }
    return $false
}

# --- Change 1: append a tab run to the "• Doanh thu theo từng suất chiếu (LICHCHIEU)" paragraph ---
$oldText1 = "• Doanh thu theo từng suất chiếu (LICHCHIEU) "
$runs1 = '<w:r><w:t xml:space="preserve">' + (Escape-Xml $oldText1) + '</w:t></w:r><w:r><w:tab/></w:r>'
Set-ParagraphRuns $oldText1 $runs1

# --- Change 2: split "7. Thống kê đồ ăn " into "7. Thống kê " + "dịch vụ" ---
$oldText2 = "7. Thống kê đồ ăn "
$newRun2a = Escape-Xml "7. Thống kê "
$newRun2b = Escape-Xml "dịch vụ"
$runs2 = '<w:r><w:t xml:space="preserve">' + $newRun2a + '</w:t></w:r><w:r><w:t>' + $newRun2b + '</w:t></w:r>'
Set-ParagraphRuns $oldText2 $runs2
